$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 181 (shifts existing rows 181-242 down to 182-243)
$ws.Rows.Item(181).Insert()

# Populate the new row 181 with the week's data (same constants as surrounding rows,
# new date + new volume/price figures)
$ws.Cells.Item(181, 1).Value = 3
$ws.Cells.Item(181, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(181, 3).Value = "Coquimbo"
$ws.Cells.Item(181, 4).Value = 44985
$ws.Cells.Item(181, 5).Value = 5
$ws.Cells.Item(181, 6).Value = 100112052
$ws.Cells.Item(181, 7).Value = "Albahaca"
$ws.Cells.Item(181, 8).Value = "Sin especificar"
$ws.Cells.Item(181, 9).Value = "Primera"
$ws.Cells.Item(181, 10).Value = 100
$ws.Cells.Item(181, 11).Value = 5000
$ws.Cells.Item(181, 12).Value = 5500
$ws.Cells.Item(181, 13).Value = 5250
$ws.Cells.Item(181, 14).Value = "$/docena de matas"
$ws.Cells.Item(181, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(181, 16).Value = 875
$ws.Cells.Item(181, 17).Value = 6
$ws.Cells.Item(181, 18).Value = "Hortaliza"

# Match the date-style formatting used by column D elsewhere (style index 2 / numFmtId 165)
$ws.Cells.Item(181, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
